$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data cells (D/E/G columns, rows 2-51) are stored as text
# (inlineStr) in the workbook, e.g. "291.47", "0.34%", "16". A plain
# Range.Value assignment of a numeric-looking string would be auto-
# coerced by Excel into a real number (losing formatting / exact text,
# and introducing floating point artifacts). Force the target range to
# Text format first so values are stored verbatim as strings, then
# restore the original (default/"Normal") style so no visual/style
# metadata is left behind on the cells.
$dataRange = $ws.Range("D2:G51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '291.65'
$ws.Range("E2").Value = '0.21%'
$ws.Range("G2").Value = '17'
$ws.Range("D3").Value = '31.08'
$ws.Range("E3").Value = '1.03%'
$ws.Range("G3").Value = '17'
$ws.Range("D4").Value = '4.952'
$ws.Range("E4").Value = '1.26%'
$ws.Range("G4").Value = '17'
$ws.Range("D5").Value = '0.07451'
$ws.Range("E5").Value = '2.68%'
$ws.Range("G5").Value = '17'
$ws.Range("D6").Value = '2.204'
$ws.Range("E6").Value = '-7.99%'
$ws.Range("G6").Value = '17'
$ws.Range("D7").Value = '7.737'
$ws.Range("E7").Value = '1.08%'
$ws.Range("G7").Value = '17'
$ws.Range("D8").Value = '0.9205'
$ws.Range("E8").Value = '2.61%'
$ws.Range("G8").Value = '17'
$ws.Range("D9").Value = '0.09416'
$ws.Range("E9").Value = '17.63%'
$ws.Range("G9").Value = '17'
$ws.Range("D10").Value = '0.1721'
$ws.Range("E10").Value = '3.16%'
$ws.Range("G10").Value = '17'
$ws.Range("D11").Value = '0.08327'
$ws.Range("E11").Value = '1.59%'
$ws.Range("G11").Value = '17'
$ws.Range("D12").Value = '0.03182'
$ws.Range("E12").Value = '3.12%'
$ws.Range("G12").Value = '17'
$ws.Range("D13").Value = '0.09928'
$ws.Range("E13").Value = '-0.87%'
$ws.Range("G13").Value = '17'
$ws.Range("D14").Value = '0.001495'
$ws.Range("E14").Value = '-0.52%'
$ws.Range("G14").Value = '17'
$ws.Range("D15").Value = '0.005688'
$ws.Range("E15").Value = '-0.74%'
$ws.Range("G15").Value = '17'
$ws.Range("E16").Value = '-0.25%'
$ws.Range("G16").Value = '17'
$ws.Range("D17").Value = '3.749'
$ws.Range("E17").Value = '1.29%'
$ws.Range("G17").Value = '17'
$ws.Range("E18").Value = '2.54%'
$ws.Range("G18").Value = '17'
$ws.Range("D19").Value = '0.3328'
$ws.Range("E19").Value = '0.29%'
$ws.Range("G19").Value = '17'
$ws.Range("E20").Value = '1.04%'
$ws.Range("G20").Value = '17'
$ws.Range("D21").Value = '4.155'
$ws.Range("E21").Value = '4.76%'
$ws.Range("G21").Value = '17'
$ws.Range("E22").Value = '-8.08%'
$ws.Range("G22").Value = '17'
$ws.Range("D23").Value = '0.04494'
$ws.Range("E23").Value = '-0.26%'
$ws.Range("G23").Value = '17'
$ws.Range("E24").Value = '0.31%'
$ws.Range("G24").Value = '17'
$ws.Range("D25").Value = '0.004260'
$ws.Range("E25").Value = '-3.45%'
$ws.Range("G25").Value = '17'
$ws.Range("D26").Value = '0.0001296'
$ws.Range("E26").Value = '-0.27%'
$ws.Range("G26").Value = '17'
$ws.Range("E27").Value = '-0.43%'
$ws.Range("G27").Value = '17'
$ws.Range("G28").Value = '17'
$ws.Range("G29").Value = '17'
$ws.Range("G30").Value = '17'
$ws.Range("G31").Value = '17'
$ws.Range("G32").Value = '17'
$ws.Range("G33").Value = '17'
$ws.Range("G34").Value = '17'
$ws.Range("G35").Value = '17'
$ws.Range("G36").Value = '17'
$ws.Range("G37").Value = '17'
$ws.Range("G38").Value = '17'
$ws.Range("D39").Value = '0.01617'
$ws.Range("E39").Value = '1.63%'
$ws.Range("G39").Value = '17'
$ws.Range("D40").Value = '0.04566'
$ws.Range("E40").Value = '4.36%'
$ws.Range("G40").Value = '17'
$ws.Range("D41").Value = '0.007404'
$ws.Range("E41").Value = '0.97%'
$ws.Range("G41").Value = '17'
$ws.Range("D42").Value = '0.009802'
$ws.Range("E42").Value = '-2.07%'
$ws.Range("G42").Value = '17'
$ws.Range("E43").Value = '3.08%'
$ws.Range("G43").Value = '17'
$ws.Range("D44").Value = '0.002154'
$ws.Range("E44").Value = '4.05%'
$ws.Range("G44").Value = '17'
$ws.Range("D45").Value = '0.009653'
$ws.Range("E45").Value = '1.68%'
$ws.Range("G45").Value = '17'
$ws.Range("D46").Value = '0.00006084'
$ws.Range("E46").Value = '6.37%'
$ws.Range("G46").Value = '17'
$ws.Range("D47").Value = '0.00000000748'
$ws.Range("E47").Value = '-0.27%'
$ws.Range("G47").Value = '17'
$ws.Range("D48").Value = '2.623'
$ws.Range("E48").Value = '16.95%'
$ws.Range("G48").Value = '17'
$ws.Range("D49").Value = '0.001991'
$ws.Range("E49").Value = '-31.30%'
$ws.Range("G49").Value = '17'
$ws.Range("D50").Value = '0.00002094'
$ws.Range("E50").Value = '-0.27%'
$ws.Range("G50").Value = '17'
$ws.Range("D51").Value = '0.0001994'
$ws.Range("E51").Value = '-0.27%'
$ws.Range("G51").Value = '17'

$dataRange.Style = "Normal"
